$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($row1, $row2) {
    $r1 = $ws.Range("B$row1`:AB$row1")
    $r2 = $ws.Range("B$row2`:AB$row2")
    $v1 = $r1.Value2
    $v2 = $r2.Value2
    $r1.Value = $v2
    $r2.Value = $v1
}

# Row 50 <-> Row 51 (match ids 7055064 / 6221723)
Swap-Rows 50 51

# Row 102 <-> Row 104 (match ids 6221814 / 6221754)
Swap-Rows 102 104

# Row 141 <-> Row 142 (match ids 8175867 / 8175866)
Swap-Rows 141 142
